$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grow the table by one row (old rows 39-41 shift down to make room); the
# exact final contents of every row from 39 down are (re)written explicitly
# below, so the shift direction here only needs to make space.
$ws.Rows(39).Insert()

# Row 39: brand-new "red.png" condition entry.
$ws.Range("A39").Value = 1
$ws.Range("B39").Value = "red.png"
$ws.Range("C39").Value = 3
$ws.Range("D39").Value = 1
$ws.Range("E39").Value = 1

# Row 40: "green.png" entry (same as the old row 40).
$ws.Range("A40").Value = 1
$ws.Range("B40").Value = "green.png"
$ws.Range("C40").Value = 3
$ws.Range("D40").Value = 0
$ws.Range("E40").Value = $null

# Row 41: "orange.png" entry (the old row 39 content, now at row 41).
$ws.Range("A41").Value = 1
$ws.Range("B41").Value = "orange.png"
$ws.Range("C41").Value = 3
$ws.Range("D41").Value = 0
$ws.Range("E41").Value = $null

# Row 42: "green.png" entry, duration corrected from 2.9 to 3.
$ws.Range("A42").Value = 1
$ws.Range("B42").Value = "green.png"
$ws.Range("C42").Value = 3
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = $null

# Add a running total of column C as a new formula cell on row 4.
$ws.Range("I4").Formula = "=SUM(C:C)"

# Restore the view state (zoom level + active selection) recorded for the
# sheet after the edit.
$excel.ActiveWindow.Zoom = 139
$ws.Range("I27").Select()
